$d = $word.ActiveDocument

$replacements = @(
    @("68×13=884", "90×38=3420"),
    @("49×14=686", "93×21=1953"),
    @("14×74=1036", "83×24=1992"),
    @("78×43=3354", "63×81=5103"),
    @("50×31=1550", "13×72=936"),
    @("68×59=4012", "64×21=1344"),
    @("29×37=1073", "89×42=3738"),
    @("31×11=341", "83×43=3569"),
    @("51×48=2448", "70×97=6790"),
    @("14×17=238", "65×44=2860"),
    @("67×30=2010", "74×24=1776"),
    @("47×92=4324", "28×75=2100"),
    @("81×92=7452", "32×84=2688"),
    @("12×50=600", "15×98=1470"),
    @("50×32=1600", "31×23=713"),
    @("59×81=4779", "81×68=5508"),
    @("46×87=4002", "49×67=3283"),
    @("71×91=6461", "24×87=2088"),
    @("40×36=1440", "39×81=3159"),
    @("13×18=234", "89×55=4895"),
    @("77×30=2310", "89×20=1780"),
    @("74×13=962", "32×60=1920"),
    @("88×89=7832", "84×30=2520"),
    @("93×76=7068", "91×21=1911"),
    @("39×90=3510", "66×98=6468")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
